$d = $word.ActiveDocument

$d.Content.Find.Execute("4 (action size)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2 (action size)", 2)
